$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 8466
$ws1.Range("F7").Value = 10523
$ws1.Range("F12").Value = 498
$ws1.Range("F22").Value = 1792
$ws1.Range("F24").Value = 530
$ws1.Range("F25").Value = 337
$ws1.Range("F27").Value = 57
$ws1.Range("F28").Value = 576
$ws1.Range("F30").Value = 1149
$ws1.Range("F33").Value = 1412
$ws1.Range("F34").Value = 434
$ws1.Range("F37").Value = 19
$ws1.Range("F39").Value = 506
$ws1.Range("F41").Value = 85
$ws1.Range("F42").Value = 277
$ws1.Range("F43").Value = 363
$ws1.Range("F45").Value = 83

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 39
$ws2.Range("F17").Value = 376

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 8466
$ws4.Range("F11").Value = 10523
$ws4.Range("F21").Value = 1792
$ws4.Range("F22").Value = 530
$ws4.Range("F23").Value = 337
$ws4.Range("F26").Value = 576
$ws4.Range("F27").Value = 39
$ws4.Range("F28").Value = 1149
$ws4.Range("F33").Value = 1412
$ws4.Range("F34").Value = 434
$ws4.Range("F39").Value = 506
$ws4.Range("F42").Value = 85
$ws4.Range("F43").Value = 277
$ws4.Range("F46").Value = 376
$ws4.Range("F47").Value = 365
$ws4.Range("F49").Value = 83
